# The uploaded workbook gained a new first row on Sheet1 that advertises
# where to find the paper's data/code (Dropbox link). Every existing row
# shifts down by one to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing content down by one row...
$ws.Rows("1").Insert()

# ...and fill the freshly inserted A1 with the new announcement text.
$ws.Range("A1").Value = "Data and Code: https://www.dropbox.com/s/b7bca82ug1112km/aaai2021_cr.rar?dl=0"
